$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace curly double quotes with straight single quotes in the en_US
# localization column (C) for the four affected rows.

$ws.Range("C3").Value = "Central Kazimierz, Kawalerielki Alliance, 'The Grand Knight Territory' Kawalerielki`n"

$ws.Range("C13").Value = "[name=`"Old Craftsman`"]  I’m not blind, which is why I clearly told you earlier 'watch out for that tree.'`n"

$ws.Range("C57").Value = "[name=`"Greatmouth Mob`"]  Today’s match is fully sponsored by the Roar Guards Company, as usual. The top ten competitors in the daily points rankings will each receive one limited edition weapon 'Rioter' provided by the Roar Guards!`n"

$ws.Range("C70").Value = "[name=`"Greatmouth Mob`"]  The reason? This competitor joining us for the first time has a very special 'name'! She may look like a cute girl, but she comes from a very famous family! `n"
